$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text (e.g. "1.150", "1.000") that must
# stay exact text rather than being auto-coerced to a number (which would drop
# trailing zeros / introduce floating point noise). Force Text format per cell
# right before writing its value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.211.60"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.800.97"
$ws.Range("E3").Value = "  +2.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.53"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4542"
$ws.Range("E7").Value = "  +20.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3539"
$ws.Range("E8").Value = "  +5.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.55"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.150"
$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07567"
$ws.Range("E11").Value = "  +4.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.74"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.003"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.250"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.272"
$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.799.55"
$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001091"
$ws.Range("E17").Value = "  +3.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06688"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.68"
$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.26"
$ws.Range("E21").Value = "  +1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.429"
$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.196.51"
$ws.Range("E23").Value = "  +1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.00"
$ws.Range("E24").Value = "  +2.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.76"
$ws.Range("E26").Value = "  +4.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.72"
$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.421"
$ws.Range("E28").Value = "  +3.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.003.90"
$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.293"
$ws.Range("E30").Value = "  -13.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.34"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  +1.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.942"
$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09478"
$ws.Range("E34").Value = "  +8.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02389"
$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.18"
$ws.Range("E36").Value = "  -1.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6743"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06295"
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2179"
$ws.Range("E39").Value = "  +3.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.204"
$ws.Range("E40").Value = "  +0.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.484"
$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.160"
$ws.Range("E43").Value = "  +1.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.03"
$ws.Range("E45").Value = "  +1.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.870"
$ws.Range("E46").Value = "  +0.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6127"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.74"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.039"
$ws.Range("E49").Value = "  +1.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07105"
$ws.Range("E50").Value = "  -2.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.171"
$ws.Range("E51").Value = "  -0.86%  "
